# Update res_bus/vm_pu.xlsx voltage-magnitude results for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.051949480937346
$ws.Range("D2").Value = 1.065145065419707
$ws.Range("E2").Value = 1.059725950246125
$ws.Range("F2").Value = 1.072610349036436
$ws.Range("I2").Value = 1.050301235574709
$ws.Range("J2").Value = 1.056974657534808
$ws.Range("K2").Value = 1.06785909540213
$ws.Range("L2").Value = 1.062454669563306
$ws.Range("M2").Value = 1.075304406058649
$ws.Range("N2").Value = 1.022744879528077
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.052983580757472
$ws.Range("D3").Value = 1.065880925245218
$ws.Range("E3").Value = 1.060613607450938
$ws.Range("F3").Value = 1.073503613925229
$ws.Range("I3").Value = 1.050569804805662
$ws.Range("J3").Value = 1.057658633487758
$ws.Range("K3").Value = 1.068410496177307
$ws.Range("L3").Value = 1.063156427115526
$ws.Range("M3").Value = 1.076014259505734
$ws.Range("N3").Value = 1.02297761036921
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053653200178864
$ws.Range("D4").Value = 1.066357405803037
$ws.Range("E4").Value = 1.061188719102321
$ws.Range("F4").Value = 1.074082320460776
$ws.Range("I4").Value = 1.050742653170084
$ws.Range("J4").Value = 1.058101089575183
$ws.Range("K4").Value = 1.068766948493543
$ws.Range("L4").Value = 1.06361060975824
$ws.Range("M4").Value = 1.076473646076237
$ws.Range("N4").Value = 1.023128038754151
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053934824405106
$ws.Range("D5").Value = 1.06655779593171
$ws.Range("E5").Value = 1.061430671650067
$ws.Range("F5").Value = 1.074325776090624
$ws.Range("I5").Value = 1.050815094566388
$ws.Range("J5").Value = 1.058287068168056
$ws.Range("K5").Value = 1.068916718501652
$ws.Range("L5").Value = 1.063801570889876
$ws.Range("M5").Value = 1.076666786495708
$ws.Range("N5").Value = 1.023191239232165
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053982117151856
$ws.Range("D6").Value = 1.066591446834442
$ws.Range("E6").Value = 1.061471306813681
$ws.Range("F6").Value = 1.074366663142739
$ws.Range("I6").Value = 1.050827244642875
$ws.Range("J6").Value = 1.058318293015289
$ws.Range("K6").Value = 1.068941860677104
$ws.Range("L6").Value = 1.063833635413423
$ws.Range("M6").Value = 1.076699216431883
$ws.Range("N6").Value = 1.023201848537442
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053656962799538
$ws.Range("D7").Value = 1.066360083122073
$ws.Range("E7").Value = 1.061191951395442
$ws.Range("F7").Value = 1.07408557287059
$ws.Range("I7").Value = 1.050743622016552
$ws.Range("J7").Value = 1.058103574748804
$ws.Range("K7").Value = 1.068768950053251
$ws.Range("L7").Value = 1.063613161301904
$ws.Range("M7").Value = 1.076476226773013
$ws.Range("N7").Value = 1.023128883398123
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.052298858512081
$ws.Range("D8").Value = 1.065393683566723
$ws.Range("E8").Value = 1.060025784673981
$ws.Range("F8").Value = 1.072912085488272
$ws.Range("I8").Value = 1.050392192880889
$ws.Range("J8").Value = 1.057205835257859
$ws.Range("K8").Value = 1.068045513884537
$ws.Range("L8").Value = 1.062691810889071
$ws.Range("M8").Value = 1.075544290518204
$ws.Range("N8").Value = 1.02282356576148
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.049909463796273
$ws.Range("D9").Value = 1.06369334968572
$ws.Range("E9").Value = 1.057976550803803
$ws.Range("F9").Value = 1.070849700429148
$ws.Range("I9").Value = 1.04976579868358
$ws.Range("J9").Value = 1.055623000147556
$ws.Range("K9").Value = 1.06676815802904
$ws.Range("L9").Value = 1.06106906632969
$ws.Range("M9").Value = 1.073902633389256
$ws.Range("N9").Value = 1.022284316015444
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048319090682833
$ws.Range("D10").Value = 1.062561609525225
$ws.Range("E10").Value = 1.056614291551687
$ws.Range("F10").Value = 1.06947851375416
$ws.Range("I10").Value = 1.049343436197162
$ws.Range("J10").Value = 1.054567210849264
$ws.Range("K10").Value = 1.065914915764914
$ws.Range("L10").Value = 1.059987819160865
$ws.Range("M10").Value = 1.072808612555321
$ws.Range("N10").Value = 1.02192400192362
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047631053589015
$ws.Range("D11").Value = 1.062072000313595
$ws.Range("E11").Value = 1.056025355578956
$ws.Range("F11").Value = 1.068885676094689
$ws.Range("I11").Value = 1.049159423842564
$ws.Range("J11").Value = 1.054109916582076
$ws.Range("K11").Value = 1.065545067414227
$ws.Range("L11").Value = 1.0595197755069
$ws.Range("M11").Value = 1.072335000454911
$ws.Range("N11").Value = 1.021767792940113
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.047375577236339
$ws.Range("D12").Value = 1.061890205292405
$ws.Range("E12").Value = 1.055806739465434
$ws.Range("F12").Value = 1.068665605463892
$ws.Range("I12").Value = 1.049090904631578
$ws.Range("J12").Value = 1.053940038141905
$ws.Range("K12").Value = 1.065407631743415
$ws.Range("L12").Value = 1.059345945587597
$ws.Range("M12").Value = 1.07215909670951
$ws.Range("N12").Value = 1.021709741647211
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04743037364555
$ws.Range("D13").Value = 1.061929197860298
$ws.Range("E13").Value = 1.055853626970228
$ws.Range("F13").Value = 1.068712805202916
$ws.Range("I13").Value = 1.049105609867862
$ws.Range("J13").Value = 1.053976478493215
$ws.Range("K13").Value = 1.065437114750208
$ws.Range("L13").Value = 1.059383231672294
$ws.Range("M13").Value = 1.072196827881384
$ws.Range("N13").Value = 1.021722195123135
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047609933970186
$ws.Range("D14").Value = 1.062056971693243
$ws.Range("E14").Value = 1.056007281822021
$ws.Range("F14").Value = 1.068867482213225
$ws.Range("I14").Value = 1.049153763464773
$ws.Range("J14").Value = 1.054095874760958
$ws.Range("K14").Value = 1.065533708107164
$ws.Range("L14").Value = 1.059505406211468
$ws.Range("M14").Value = 1.072320459846788
$ws.Range("N14").Value = 1.021762994978054
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04772057911632
$ws.Range("D15").Value = 1.062135706369781
$ws.Range("E15").Value = 1.056101972352358
$ws.Range("F15").Value = 1.068962801837306
$ws.Range("I15").Value = 1.04918341012652
$ws.Range("J15").Value = 1.054169436244647
$ws.Range("K15").Value = 1.065593214866273
$ws.Range("L15").Value = 1.059580684959892
$ws.Range("M15").Value = 1.07239663584516
$ws.Range("N15").Value = 1.021788129371825
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048364766216019
$ws.Range("D16").Value = 1.062594112691827
$ws.Range("E16").Value = 1.05665339700011
$ws.Range("F16").Value = 1.069517877398703
$ws.Range("I16").Value = 1.049355624791456
$ws.Range("J16").Value = 1.054597557266613
$ws.Range("K16").Value = 1.06593945326583
$ws.Range("L16").Value = 1.060018884766167
$ws.Range("M16").Value = 1.072840046946207
$ws.Range("N16").Value = 1.021934365010454
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048769009988765
$ws.Range("D17").Value = 1.062881778033217
$ws.Range("E17").Value = 1.056999541158105
$ws.Range("F17").Value = 1.06986630213855
$ws.Range("I17").Value = 1.049463349178802
$ws.Range("J17").Value = 1.054866071764068
$ws.Range("K17").Value = 1.066156535996154
$ws.Range("L17").Value = 1.060293794893313
$ws.Range("M17").Value = 1.07311821608782
$ws.Range("N17").Value = 1.022026043955145
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049004856886725
$ws.Range("D18").Value = 1.063049610896742
$ws.Range("E18").Value = 1.057201530975296
$ws.Range("F18").Value = 1.070069618691214
$ws.Range("I18").Value = 1.04952607426072
$ws.Range("J18").Value = 1.055022679090956
$ws.Range("K18").Value = 1.066283119104788
$ws.Range("L18").Value = 1.060454159113563
$ws.Range("M18").Value = 1.073280477620402
$ws.Range("N18").Value = 1.02207950032084
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049085284469176
$ws.Range("D19").Value = 1.06310684474227
$ws.Range("E19").Value = 1.057270419481474
$ws.Range("F19").Value = 1.070138958983188
$ws.Range("I19").Value = 1.049547443437005
$ws.Range("J19").Value = 1.05507607598004
$ws.Range("K19").Value = 1.066326274273167
$ws.Range("L19").Value = 1.060508841484723
$ws.Range("M19").Value = 1.073335806307987
$ws.Range("N19").Value = 1.022097724427594
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048725632420076
$ws.Range("D20").Value = 1.062850909855531
$ws.Range("E20").Value = 1.056962393871453
$ws.Range("F20").Value = 1.06982891052682
$ws.Range("I20").Value = 1.049451802614928
$ws.Range("J20").Value = 1.054837264002207
$ws.Range("K20").Value = 1.066133248953771
$ws.Range("L20").Value = 1.060264298205974
$ws.Range("M20").Value = 1.073088370111931
$ws.Range("N20").Value = 1.022016209574138
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047557055407804
$ws.Range("D21").Value = 1.062019343606713
$ws.Range("E21").Value = 1.055962030399414
$ws.Range("F21").Value = 1.068821929934255
$ws.Range("I21").Value = 1.049139588087719
$ws.Range("J21").Value = 1.054060716055133
$ws.Range("K21").Value = 1.065505265351909
$ws.Range("L21").Value = 1.059469428232068
$ws.Range("M21").Value = 1.072284052841078
$ws.Range("N21").Value = 1.021750981216963
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.046822852843425
$ws.Range("D22").Value = 1.061496896874225
$ws.Range("E22").Value = 1.055333878056871
$ws.Range("F22").Value = 1.068189586808088
$ws.Range("I22").Value = 1.048942309843654
$ws.Range("J22").Value = 1.053572359989998
$ws.Range("K22").Value = 1.065110094751336
$ws.Range("L22").Value = 1.058969791530006
$ws.Range("M22").Value = 1.07177844475219
$ws.Range("N22").Value = 1.021584057920242
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.047212017323113
$ws.Range("D23").Value = 1.061773818192872
$ws.Range("E23").Value = 1.055666795823344
$ws.Range("F23").Value = 1.068524729010722
$ws.Range("I23").Value = 1.049046983197182
$ws.Range("J23").Value = 1.053831256955321
$ws.Range("K23").Value = 1.065319613372404
$ws.Range("L23").Value = 1.059234645855906
$ws.Range("M23").Value = 1.072046467556568
$ws.Range("N23").Value = 1.021672562529912
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048745232695873
$ws.Range("D24").Value = 1.062864857725018
$ws.Range("E24").Value = 1.056979178853956
$ws.Range("F24").Value = 1.069845805920116
$ws.Range("I24").Value = 1.049457020346112
$ws.Range("J24").Value = 1.05485028102884
$ws.Range("K24").Value = 1.066143771482021
$ws.Range("L24").Value = 1.060277626446735
$ws.Range("M24").Value = 1.073101856192349
$ws.Range("N24").Value = 1.02202065336451
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050526731162515
$ws.Range("D25").Value = 1.064132612371233
$ws.Range("E25").Value = 1.058505644735652
$ws.Range("F25").Value = 1.071382223107872
$ws.Range("I25").Value = 1.049928578777183
$ws.Range("J25").Value = 1.056032303756731
$ws.Range("K25").Value = 1.06709868384963
$ws.Range("L25").Value = 1.06106906632969
$ws.Range("M25").Value = 1.074326971850839
$ws.Range("N25").Value = 1.02242386999661

Write-Output "Updated 264 cells"
